$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "cliente" values that were removed for rows 6, 8 and 10
$ws.Range("D6").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("D10").ClearContents()

# Update the active selection to match the saved cursor position
$ws.Range("D8").Select()
